$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5877678990364075
$ws.Range("B1").Value = 1.240883111953735
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.591520309448242
$ws.Range("E1").Value = 1.421040534973145
